$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header updates: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 32   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/5/2025  Through  5/11/2025"

# --- Crime Complaints table updates (rows 15-30) ---
$ws.Range("C15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 2
$ws.Range("F14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("F14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = -50
$ws.Range("K14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 166.666666666667
$ws.Range("C16").Value = 4
$ws.Range("F14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = "'***.*"
$ws.Range("C14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 6
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 36
$ws.Range("K16").Value = -35.714285714285
$ws.Range("L16").Value = -32.075471698113
$ws.Range("M16").Value = -32.075471698113
$ws.Range("N16").Value = -84.279475982532
$ws.Range("C17").Value = 6
$ws.Range("F14").Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = -10.526315789473
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 69
$ws.Range("K17").Value = 15.942028985507
$ws.Range("L17").Value = 11.111111111111
$ws.Range("M17").Value = 95.121951219512
$ws.Range("N17").Value = -10.112359550561
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -10
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = 22.222222222222
$ws.Range("L18").Value = -19.117647058823
$ws.Range("M18").Value = 22.222222222222
$ws.Range("N18").Value = -86.352357320099
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 27.777777777777
$ws.Range("I19").Value = 210
$ws.Range("J19").Value = 269
$ws.Range("K19").Value = -21.933085501858
$ws.Range("L19").Value = -9.090909090909
$ws.Range("M19").Value = 87.5
$ws.Range("N19").Value = 39.072847682119
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 21
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 23.529411764705
$ws.Range("I20").Value = 87
$ws.Range("J20").Value = 124
$ws.Range("K20").Value = -29.838709677419
$ws.Range("L20").Value = -23.684210526315
$ws.Range("M20").Value = 93.333333333333
$ws.Range("N20").Value = -86.363636363636
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = 14.285714285714
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = 10.714285714285
$ws.Range("I21").Value = 477
$ws.Range("J21").Value = 567
$ws.Range("K21").Value = -15.873015873015
$ws.Range("L21").Value = -12.637362637362
$ws.Range("M21").Value = 57.947019867549
$ws.Range("N21").Value = -68.618421052631
$ws.Range("C22").Value = 1
$ws.Range("F14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("F14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 100
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = -50
$ws.Range("D23").Value = 2
$ws.Range("F14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("F23").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 26
$ws.Range("K23").Value = -38.461538461538
$ws.Range("L23").Value = -40.740740740740
$ws.Range("M23").Value = -5.882352941176
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 13
$ws.Range("E24").Value = -30.769230769230
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = 32.876712328767
$ws.Range("I24").Value = 352
$ws.Range("J24").Value = 391
$ws.Range("K24").Value = -9.974424552429
$ws.Range("L24").Value = -7.853403141361
$ws.Range("M24").Value = 16.943521594684
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 38
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 15.151515151515
$ws.Range("I25").Value = 103
$ws.Range("J25").Value = 164
$ws.Range("K25").Value = -37.195121951219
$ws.Range("L25").Value = -34.394904458598
$ws.Range("C26").Value = 15
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 47
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = 80.769230769230
$ws.Range("I26").Value = 145
$ws.Range("J26").Value = 117
$ws.Range("K26").Value = 23.931623931623
$ws.Range("L26").Value = 25
$ws.Range("M26").Value = 16.935483870967
$ws.Range("D27").Value = 2
$ws.Range("F14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -50
$ws.Range("K14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -33.333333333333
$ws.Range("I27").Value = 10
$ws.Range("J27").Value = 8
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = -16.666666666666
$ws.Range("C28").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = 4
$ws.Range("F14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 13
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -45.833333333333
$ws.Range("L28").Value = 44.444444444444
$ws.Range("F29").Value = 3
$ws.Range("F30").Value = 2
